$d = $word.ActiveDocument

# --- Edit 1: "Oligodendrocyte derived from oligodendrocyte progenitor cells (OPCs) produce"
#             -> "Oligodendrocytes derived from oligodendrocyte progenitor cells (OPCs) produce"
$d.Content.Find.Execute("Oligodendrocyte derived from oligodendrocyte progenitor cells (OPCs) produce",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Oligodendrocytes derived from oligodendrocyte progenitor cells (OPCs) produce", 2)

# Split the newly-inserted "s" into its own run (matching the authored edit, which left the
# pluralizing "s" as a distinct run) by toggling a character property on/off.
$rng = $d.Content
$rng.Find.Execute("Oligodendrocytes derived", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$sRange = $d.Range($rng.Start + 15, $rng.Start + 16)
$sRange.Bold = 1
$sRange.Bold = 0

# --- Edit 2: "promoting to the production" -> "promoting the production"
$d.Content.Find.Execute("promoting to the production", $true, $false, $false, $false, $false,
                         $true, 1, $false, "promoting the production", 2)
